$wb = $excel.ActiveWorkbook

# ----- Sheet: ALC -----
$ws = $wb.Worksheets.Item("ALC")
# Row 9
$ws.Range("H9").Value = 112.333336
$ws.Range("I9").Value = 133.2
$ws.Range("J9").Value = 8
$ws.Range("K9").Value = 133.2
$ws.Range("L9").Value = 8
$ws.Range("M9").Value = 35.80000000000001
$ws.Range("N9").Value = -346
# Row 51
$ws.Range("H51").Value = 5399
$ws.Range("I51").Value = 2995
$ws.Range("K51").Value = 2995
$ws.Range("M51").Value = -2511
# Row 70
$ws.Range("H70").Value = 9999.666999999999
$ws.Range("J70").Value = 9999.5
$ws.Range("L70").Value = 29998.5
$ws.Range("N70").Value = -30538.5
# Row 73
$ws.Range("H73").Value = 9999.666999999999
$ws.Range("J73").Value = 9999.5
$ws.Range("L73").Value = 29998.5
$ws.Range("N73").Value = -31870.5
# Row 94
$ws.Range("H94").Value = 3035.1428
$ws.Range("I94").Value = 1957.6666
$ws.Range("K94").Value = 1957.6666
$ws.Range("M94").Value = -1506.6666
# Row 95
$ws.Range("H95").Value = 27497.5
$ws.Range("J95").Value = 27497.5
$ws.Range("L95").Value = 27497.5
$ws.Range("N95").Value = -32989.5
# Row 100
$ws.Range("H100").Value = 0
$ws.Range("I100").Value = 0
$ws.Range("K100").Value = 0
$ws.Range("M100").ClearContents()

# ----- Sheet: ARM -----
$ws = $wb.Worksheets.Item("ARM")
# Row 101
$ws.Range("H101").Value = 65995.336
$ws.Range("J101").Value = 65995.336
$ws.Range("L101").Value = 65995.336
$ws.Range("N101").Value = -72485.336
# Row 110
$ws.Range("H110").Value = 111114264
$ws.Range("I110").Value = 200002500
$ws.Range("K110").Value = 200002500
$ws.Range("M110").Value = -200000455

# ----- Sheet: BSM -----
$ws = $wb.Worksheets.Item("BSM")
# Row 86
$ws.Range("H86").Value = 4737.4707
$ws.Range("I86").Value = 2753.0833
$ws.Range("K86").Value = 2753.0833
$ws.Range("M86").Value = -1630.0833
# Row 89
$ws.Range("H89").Value = 4737.4707
$ws.Range("I89").Value = 2753.0833
$ws.Range("K89").Value = 13765.4165
$ws.Range("M89").Value = -8149.416499999999

# ----- Sheet: CRP -----
$ws = $wb.Worksheets.Item("CRP")
# Row 16
$ws.Range("H16").Value = 2396.5
$ws.Range("I16").Value = 1528.8334
$ws.Range("J16").Value = 4999.5
$ws.Range("K16").Value = 1528.8334
$ws.Range("L16").Value = 4999.5
$ws.Range("M16").Value = -1241.8334
$ws.Range("N16").Value = -5573.5
# Row 19
$ws.Range("H19").Value = 390.81818
$ws.Range("I19").Value = 349.8
$ws.Range("K19").Value = 349.8
$ws.Range("M19").Value = -179.8
# Row 24
$ws.Range("H24").Value = 390.81818
$ws.Range("I24").Value = 349.8
$ws.Range("K24").Value = 349.8
$ws.Range("M24").Value = -179.8
# Row 113
$ws.Range("H113").Value = 2396.5
$ws.Range("I113").Value = 1528.8334
$ws.Range("J113").Value = 4999.5
$ws.Range("K113").Value = 1528.8334
$ws.Range("L113").Value = 4999.5
$ws.Range("M113").Value = 641.1666
$ws.Range("N113").Value = -9339.5

# ----- Sheet: CUL -----
$ws = $wb.Worksheets.Item("CUL")
# Row 60
$ws.Range("H60").Value = 947.6842
$ws.Range("I60").Value = 214.71428
$ws.Range("J60").Value = 3000
$ws.Range("K60").Value = 644.14284
$ws.Range("L60").Value = 9000
$ws.Range("M60").Value = -393.14284
$ws.Range("N60").Value = -9502
# Row 92
$ws.Range("H92").Value = 4170
$ws.Range("I92").Value = 2000
$ws.Range("J92").Value = 4712.5
$ws.Range("K92").Value = 6000
$ws.Range("L92").Value = 14137.5
$ws.Range("M92").Value = -4752
$ws.Range("N92").Value = -16633.5
# Row 104
$ws.Range("H104").Value = 9180
$ws.Range("I104").Value = 6250
$ws.Range("K104").Value = 18750
$ws.Range("M104").Value = -16129
# Row 129
$ws.Range("H129").Value = 2266.4
$ws.Range("I129").Value = 720
$ws.Range("J129").Value = 3812.8
$ws.Range("K129").Value = 2160
$ws.Range("L129").Value = 11438.4
$ws.Range("M129").Value = 2840
$ws.Range("N129").Value = -21438.4
# Row 130
$ws.Range("H130").Value = 2566.6667
$ws.Range("I130").Value = 2566.6667
$ws.Range("K130").Value = 7700.000100000001
$ws.Range("M130").Value = -2680.000100000001
# Row 131
$ws.Range("H131").Value = 1847.0526
$ws.Range("J131").Value = 2810.7778
$ws.Range("L131").Value = 8432.3334
$ws.Range("N131").Value = -18512.3334
# Row 134
$ws.Range("H134").Value = 2350
$ws.Range("I134").Value = 2350
$ws.Range("K134").Value = 7050
$ws.Range("M134").Value = -1980
# Row 138
$ws.Range("H138").Value = 8944.444
$ws.Range("I138").Value = 5000
$ws.Range("J138").Value = 9437.5
$ws.Range("K138").Value = 15000
$ws.Range("L138").Value = 28312.5
$ws.Range("M138").Value = -9860
$ws.Range("N138").Value = -38592.5
# Row 139
$ws.Range("H139").Value = 1057.6
$ws.Range("I139").Value = 1057.6
$ws.Range("J139").Value = 0
$ws.Range("K139").Value = 3172.8
$ws.Range("L139").Value = 0
$ws.Range("M139").ClearContents()
$ws.Range("N139").Value = 1967.2
# Row 140
$ws.Range("H140").Value = 2718.6
$ws.Range("I140").Value = 2125.375
$ws.Range("J140").Value = 5091.5
$ws.Range("K140").Value = 6376.125
$ws.Range("L140").Value = 15274.5
$ws.Range("M140").Value = -1196.125
$ws.Range("N140").Value = -25634.5
# Row 141
$ws.Range("H141").Value = 1749
$ws.Range("I141").Value = 1749
$ws.Range("K141").Value = 5247
$ws.Range("M141").Value = -67

# ----- Sheet: GSM -----
$ws = $wb.Worksheets.Item("GSM")
# Row 14
$ws.Range("H14").Value = 549.5
$ws.Range("J14").Value = 549.5
$ws.Range("L14").Value = 549.5
$ws.Range("N14").Value = -885.5
# Row 102
$ws.Range("H102").Value = 857.3333
$ws.Range("I102").Value = 978.8
$ws.Range("J102").Value = 250
$ws.Range("K102").Value = 978.8
$ws.Range("L102").Value = 250
$ws.Range("M102").Value = 643.2
$ws.Range("N102").Value = -3494

# ----- Sheet: LTW -----
$ws = $wb.Worksheets.Item("LTW")
# Row 61
$ws.Range("H61").Value = 125005910
$ws.Range("I61").Value = 200004660
$ws.Range("K61").Value = 200004660
$ws.Range("M61").Value = -200004458
# Row 100
$ws.Range("H100").Value = 4712.727
$ws.Range("I100").Value = 1068
$ws.Range("J100").Value = 7750
$ws.Range("K100").Value = 1068
$ws.Range("L100").Value = 7750
$ws.Range("M100").Value = -527
$ws.Range("N100").Value = -8832
# Row 113
$ws.Range("H113").Value = 125005910
$ws.Range("I113").Value = 200004660
$ws.Range("K113").Value = 200004660
$ws.Range("M113").Value = -200002490

# ----- Sheet: WVR -----
$ws = $wb.Worksheets.Item("WVR")
# Row 6
$ws.Range("H6").Value = 566.6667
$ws.Range("J6").Value = 600
$ws.Range("L6").Value = 600
$ws.Range("N6").Value = -830
# Row 11
$ws.Range("H11").Value = 250002.5
$ws.Range("J11").Value = 250002.5
$ws.Range("L11").Value = 250002.5
$ws.Range("N11").Value = -250286.5
# Row 29
$ws.Range("H29").Value = 4066.6667
$ws.Range("J29").Value = 4066.6667
$ws.Range("L29").Value = 4066.6667
$ws.Range("M29").Value = -4646.6667
# Row 100
$ws.Range("H100").Value = 1094.4546
$ws.Range("I100").Value = 1328.4286
$ws.Range("J100").Value = 685
$ws.Range("K100").Value = 2656.8572
$ws.Range("L100").Value = 1370
$ws.Range("M100").Value = -2115.8572
$ws.Range("N100").Value = -2452
# Row 132
$ws.Range("H132").Value = 929.3158
$ws.Range("I132").Value = 950.6667
$ws.Range("K132").Value = 2852.0001
$ws.Range("M132").Value = -322.0001000000002
